$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the "Valid Login" test script
$ws.Name = "ValidLogin"

# Populate the login-test fixture data. Values are written in the same
# order the shared-string table ends up in (admin/pointofsale first,
# then the header labels) so the data reflects how the sheet was built.
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pointofsale"
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

# Column B holds the longer "pointofsale" value - autofit it so it is
# fully visible.
$ws.Columns("B").AutoFit() | Out-Null

# Zoom in on the sheet and move the selection down to A3, below the data.
$excel.ActiveWindow.Zoom = 205
$ws.Range("A3").Select() | Out-Null
